# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.232.53'
$ws.Range('E2').Value = '  -0.69%  '
$ws.Range('D3').Value = '1.656.10'
$ws.Range('E3').Value = '  -1.01%  '
$ws.Range('E4').Value = '  -0.62%  '
$ws.Range('D5').Value = '''219.28'
$ws.Range('E5').Value = '  -0.73%  '
$ws.Range('D6').Value = '''0.5229'
$ws.Range('E6').Value = '  -2.53%  '
$ws.Range('D8').Value = '''0.2657'
$ws.Range('E8').Value = '  -0.67%  '
$ws.Range('D9').Value = '''0.06328'
$ws.Range('E9').Value = '  -1.42%  '
$ws.Range('D10').Value = '''20.64'
$ws.Range('E10').Value = '  -1.89%  '
$ws.Range('D11').Value = '''0.07775'
$ws.Range('E11').Value = '  -0.95%  '
$ws.Range('D12').Value = '''4.558'
$ws.Range('E12').Value = '  -0.26%  '
$ws.Range('D13').Value = '1.655.77'
$ws.Range('E13').Value = '  -1.07%  '
$ws.Range('D14').Value = '1.885.49'
$ws.Range('E14').Value = '  -0.92%  '
$ws.Range('D15').Value = '''0.5657'
$ws.Range('E15').Value = '  +0.10%  '
$ws.Range('D16').Value = '0.0₅8114'
$ws.Range('E16').Value = '  -1.12%  '
$ws.Range('D17').Value = '''65.43'
$ws.Range('E17').Value = '  -1.64%  '
$ws.Range('D18').Value = '26.233.03'
$ws.Range('E18').Value = '  -0.82%  '
$ws.Range('E20').Value = '  +0.14%  '
$ws.Range('D21').Value = '''193.02'
$ws.Range('E21').Value = '  -2.51%  '
$ws.Range('E22').Value = '  -0.63%  '
$ws.Range('D23').Value = '''6.030'
$ws.Range('E23').Value = '  -0.83%  '
$ws.Range('E24').Value = '  -0.68%  '
$ws.Range('D25').Value = '''143.94'
$ws.Range('E25').Value = '  -1.81%  '
$ws.Range('D26').Value = '''0.1202'
$ws.Range('E26').Value = '  -2.68%  '
$ws.Range('D27').Value = '''7.269'
$ws.Range('E27').Value = '  -0.11%  '
$ws.Range('D28').Value = '''16.00'
$ws.Range('E28').Value = '  -1.76%  '
$ws.Range('D29').Value = '''1.499'
$ws.Range('E29').Value = '  -1.03%  '
$ws.Range('D30').Value = '''0.05610'
$ws.Range('E30').Value = '  -4.87%  '
$ws.Range('D31').Value = '''1.278'
$ws.Range('E31').Value = '  -1.04%  '
$ws.Range('D32').Value = '''3.502'
$ws.Range('E32').Value = '  -2.33%  '
$ws.Range('D33').Value = '''3.382'
$ws.Range('E33').Value = '  +2.09%  '
$ws.Range('D34').Value = '''1.589'
$ws.Range('E34').Value = '  -2.25%  '
$ws.Range('D35').Value = '''2.804'
$ws.Range('E35').Value = '  -1.87%  '
$ws.Range('D36').Value = '''0.9447'
$ws.Range('E36').Value = '  -2.86%  '
$ws.Range('D37').Value = '''2.406'
$ws.Range('E37').Value = '  -1.09%  '
$ws.Range('D38').Value = '''0.5758'
$ws.Range('E38').Value = '  -1.44%  '
$ws.Range('D39').Value = '''0.01601'
$ws.Range('E39').Value = '  -0.85%  '
$ws.Range('D40').Value = '''5.920'
$ws.Range('E40').Value = '  -0.05%  '
$ws.Range('D41').Value = '''2.584'
$ws.Range('E41').Value = '  -0.45%  '
$ws.Range('D42').Value = '''0.8478'
$ws.Range('E42').Value = '  -2.20%  '
$ws.Range('E43').Value = '  -0.70%  '
$ws.Range('D44').Value = '1.036.51'
$ws.Range('E44').Value = '  -4.01%  '
$ws.Range('D45').Value = '''102.39'
$ws.Range('E45').Value = '  -2.02%  '
$ws.Range('D46').Value = '1.795.86'
$ws.Range('E46').Value = '  -0.91%  '
$ws.Range('D47').Value = '''58.48'
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').Value = '0.0₈104'
$ws.Range('E48').Value = '  -2.17%  '
$ws.Range('D49').Value = '''1.003'
$ws.Range('E49').Value = '  -1.16%  '
$ws.Range('D50').Value = '''0.05316'
$ws.Range('E50').Value = '  +2.88%  '
$ws.Range('E51').Value = '  -1.21%  '
